$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Updated Gaussian Filter benchmark numbers ---
$ws.Range("B15").Value = 0.267
$ws.Range("F15").Value = 11.948
$ws.Range("F16").Value = 19.204

# --- Remove the old "Why is this not much slower than NOOP?" note on J16 ---
$ws.Range("J16").ClearContents() | Out-Null

# --- New row 17: full implementation (+ GetOffsetImageElement) ---
$ws.Range("C17").Value = "full implementation (+ GetOffsetImageElement)"
$ws.Range("F17").Value = 103.049
$ws.Range("H17").Formula = "=F17/`$B`$15"

# --- New "512 Depth Buffer" mini benchmark block (rows 6-7, cols L:N) ---
$ws.Range("L6").Value = "512 Depth Buffer"
$ws.Range("L7").Value = 6.066
$ws.Range("M7").Value = 0.011
$ws.Range("N7").Formula = "=L7/M7"

# --- Restore the active-cell selection to match the authored state ---
$ws.Range("L6").Select() | Out-Null
